# Fruta / hortaliza, semanal
# Insert the latest week's two rows (Primera / Segunda quality) for
# Brocoli @ Terminal La Palmera de La Serena, shifting all the
# existing historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current row 726 (the first
# data row of the table's body at that point); this pushes the
# existing rows 726:780 down to 728:782 and copies formatting
# (including the date-style D column) from the row above.
$ws.Rows("726:727").Insert()

# --- Row 726: "Primera" quality, week of 2022-08-10 (serial 44783) ---
$ws.Cells.Item(726, 1).Value = 8
$ws.Cells.Item(726, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(726, 3).Value = "Coquimbo"
$ws.Cells.Item(726, 4).Value = 44783
$ws.Cells.Item(726, 5).Value = 4
$ws.Cells.Item(726, 6).Value = 100112023
$ws.Cells.Item(726, 7).Value = "Brócoli"
$ws.Cells.Item(726, 8).Value = "Sin especificar"
$ws.Cells.Item(726, 9).Value = "Primera"
$ws.Cells.Item(726, 10).Value = 2600
$ws.Cells.Item(726, 11).Value = 750
$ws.Cells.Item(726, 12).Value = 800
$ws.Cells.Item(726, 13).Value = 775
$ws.Cells.Item(726, 14).Value = "`$/unidad"
$ws.Cells.Item(726, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(726, 16).Value = 775
$ws.Cells.Item(726, 17).Value = 1
$ws.Cells.Item(726, 18).Value = "Hortaliza"

# --- Row 727: "Segunda" quality, week of 2022-08-10 (serial 44783) ---
$ws.Cells.Item(727, 1).Value = 8
$ws.Cells.Item(727, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(727, 3).Value = "Coquimbo"
$ws.Cells.Item(727, 4).Value = 44783
$ws.Cells.Item(727, 5).Value = 4
$ws.Cells.Item(727, 6).Value = 100112023
$ws.Cells.Item(727, 7).Value = "Brócoli"
$ws.Cells.Item(727, 8).Value = "Sin especificar"
$ws.Cells.Item(727, 9).Value = "Segunda"
$ws.Cells.Item(727, 10).Value = 1500
$ws.Cells.Item(727, 11).Value = 650
$ws.Cells.Item(727, 12).Value = 700
$ws.Cells.Item(727, 13).Value = 675
$ws.Cells.Item(727, 14).Value = "`$/unidad"
$ws.Cells.Item(727, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(727, 16).Value = 675
$ws.Cells.Item(727, 17).Value = 1
$ws.Cells.Item(727, 18).Value = "Hortaliza"
